# Fruta / hortaliza, semanal
# Swap the weekly price data between the two reporting dates:
#   rows 2-3 (date 44216) <-> rows 4-5 (date 44195)
# Columns touched: D (Fecha), N (Precio mínimo), O (Precio máximo),
#                  P (Precio promedio ponderado), S (Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap row 2 <-> row 4 ---
$D2 = $ws.Range("D2").Value2
$N2 = $ws.Range("N2").Value2
$O2 = $ws.Range("O2").Value2
$P2 = $ws.Range("P2").Value2
$S2 = $ws.Range("S2").Value2

$D4 = $ws.Range("D4").Value2
$N4 = $ws.Range("N4").Value2
$O4 = $ws.Range("O4").Value2
$P4 = $ws.Range("P4").Value2
$S4 = $ws.Range("S4").Value2

$ws.Range("D2").Value2 = $D4
$ws.Range("N2").Value2 = $N4
$ws.Range("O2").Value2 = $O4
$ws.Range("P2").Value2 = $P4
$ws.Range("S2").Value2 = $S4

$ws.Range("D4").Value2 = $D2
$ws.Range("N4").Value2 = $N2
$ws.Range("O4").Value2 = $O2
$ws.Range("P4").Value2 = $P2
$ws.Range("S4").Value2 = $S2

# --- swap row 3 <-> row 5 ---
$D3 = $ws.Range("D3").Value2
$N3 = $ws.Range("N3").Value2
$O3 = $ws.Range("O3").Value2
$P3 = $ws.Range("P3").Value2
$S3 = $ws.Range("S3").Value2

$D5 = $ws.Range("D5").Value2
$N5 = $ws.Range("N5").Value2
$O5 = $ws.Range("O5").Value2
$P5 = $ws.Range("P5").Value2
$S5 = $ws.Range("S5").Value2

$ws.Range("D3").Value2 = $D5
$ws.Range("N3").Value2 = $N5
$ws.Range("O3").Value2 = $O5
$ws.Range("P3").Value2 = $P5
$ws.Range("S3").Value2 = $S5

$ws.Range("D5").Value2 = $D3
$ws.Range("N5").Value2 = $N3
$ws.Range("O5").Value2 = $O3
$ws.Range("P5").Value2 = $P3
$ws.Range("S5").Value2 = $S3
